$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.47"
$ws.Range("E2").Value = "'0.31%"
$ws.Range("D3").Value = "'40.94"
$ws.Range("E3").Value = "'3.80%"
$ws.Range("E4").Value = "'0.31%"
$ws.Range("E5").Value = "'-0.76%"
$ws.Range("D6").Value = "'1.603"
$ws.Range("E6").Value = "'0.37%"
$ws.Range("D8").Value = "'0.9019"
$ws.Range("E8").Value = "'0.31%"
$ws.Range("D9").Value = "'0.1117"
$ws.Range("E9").Value = "'10.52%"
$ws.Range("D10").Value = "'0.1795"
$ws.Range("E10").Value = "'3.33%"
$ws.Range("D11").Value = "'0.09175"
$ws.Range("E11").Value = "'1.67%"
$ws.Range("D12").Value = "'0.04167"
$ws.Range("E12").Value = "'-6.61%"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("E13").Value = "'-0.12%"
$ws.Range("D14").Value = "'0.001260"
$ws.Range("E14").Value = "'-1.09%"
$ws.Range("D15").Value = "'0.005870"
$ws.Range("E15").Value = "'0.55%"
$ws.Range("D16").Value = "'3.341"
$ws.Range("E16").Value = "'-0.38%"
$ws.Range("D17").Value = "'4.246"
$ws.Range("E17").Value = "'0.29%"
$ws.Range("E18").Value = "'-0.02%"
$ws.Range("D19").Value = "'6.629"
$ws.Range("E19").Value = "'-6.09%"
$ws.Range("E20").Value = "'1.11%"
$ws.Range("E21").Value = "'-0.42%"
$ws.Range("D22").Value = "'0.04073"
$ws.Range("E22").Value = "'-1.97%"
$ws.Range("D23").Value = "'0.001246"
$ws.Range("E23").Value = "'3.04%"
$ws.Range("D24").Value = "'0.004089"
$ws.Range("E24").Value = "'0.94%"
$ws.Range("E25").Value = "'-0.18%"
$ws.Range("D38").Value = "'0.02401"
$ws.Range("E38").Value = "'3.10%"
$ws.Range("D39").Value = "'0.05203"
$ws.Range("E39").Value = "'0.99%"
$ws.Range("D40").Value = "'0.007770"
$ws.Range("E40").Value = "'-1.69%"
$ws.Range("E41").Value = "'-1.53%"
$ws.Range("D42").Value = "'0.007053"
$ws.Range("E42").Value = "'11.97%"
$ws.Range("D43").Value = "'0.001951"
$ws.Range("E43").Value = "'-0.17%"
$ws.Range("D44").Value = "'0.007710"
$ws.Range("E44").Value = "'-6.20%"
$ws.Range("E45").Value = "'-7.69%"
$ws.Range("D46").Value = "'0.00006973"
$ws.Range("E46").Value = "'7.03%"
$ws.Range("E47").Value = "'-0.18%"
$ws.Range("D48").Value = "'0.04913"
$ws.Range("E48").Value = "'1,286.72%"
$ws.Range("E50").Value = "'-0.18%"
$ws.Range("E51").Value = "'-0.18%"
